$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.849.86"
$ws.Range("E2").Value = "  +4.35%  "

$ws.Range("D3").Value = "2.261.85"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.89"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.67"
$ws.Range("E6").Value = "  +6.93%  "

$ws.Range("E7").Value = "  -1.03%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.41"
$ws.Range("E10").Value = "  +4.10%  "

$ws.Range("E11").Value = "  -2.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.09"
$ws.Range("E12").Value = "  -1.44%  "

$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").Value = "2.605.66"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").Value = "2.258.37"
$ws.Range("E15").Value = "  -0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.55"
$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").Value = "46.868.29"
$ws.Range("E17").Value = "  +4.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.790"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -5.69%  "

$ws.Range("E20").Value = "  +1.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.80"
$ws.Range("E21").Value = "  -3.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.25"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.91"
$ws.Range("E23").Value = "  +3.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  -2.08%  "

$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$ws.Range("E26").Value = "  -1.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.11"
$ws.Range("E27").Value = "  +2.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  +1.44%  "

$ws.Range("E30").Value = "  +2.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.82"
$ws.Range("E31").Value = "  +10.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "145.57"
$ws.Range("E32").Value = "  -4.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.36"
$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.23"
$ws.Range("E34").Value = "  +11.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0765"
$ws.Range("E35").Value = "  -2.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.115"
$ws.Range("E36").Value = "  +11.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("E37").Value = "  -1.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.08"
$ws.Range("E38").Value = "  +18.50%  "

$ws.Range("E39").Value = "  -4.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.87"
$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0297"
$ws.Range("E41").Value = "  -3.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("E42").Value = "  -2.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.18"
$ws.Range("E45").Value = "  +19.49%  "

$ws.Range("D46").Value = "1.772.94"
$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "71.10"
$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.185"
$ws.Range("E48").Value = "  -3.87%  "

$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.485.10"
$ws.Range("E51").Value = "  -0.14%  "
